$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 30,8
$data[0,0] = 0
$data[0,1] = "walkingToRunning"
$data[0,2] = -0.7247905336188185
$data[0,3] = -24.30215097461226
$data[0,4] = -0.4334598116620501
$data[0,5] = -1.109872460365295
$data[0,6] = 1.229648113250732
$data[0,7] = 2.535529136657715
$data[1,0] = 100
$data[1,1] = "walkingToRunning"
$data[1,2] = 0.9856911530861703
$data[1,3] = -27.7952582342386
$data[1,4] = -10.58349611914382
$data[1,5] = -0.0105194868519902
$data[1,6] = 0.2277668565511703
$data[1,7] = 1.795702934265137
$data[2,0] = 200
$data[2,1] = "walkingToRunning"
$data[2,2] = -1.217567012860227
$data[2,3] = -15.8785987887862
$data[2,4] = 2.25199632955022
$data[2,5] = -0.8973521590232849
$data[2,6] = 0.1004677563905716
$data[2,7] = 1.17305588722229
$data[3,0] = 300
$data[3,1] = "walkingToRunning"
$data[3,2] = 0.9202353192504198
$data[3,3] = -13.56062749151649
$data[3,4] = -2.717396190180546
$data[3,5] = -0.6182528734207153
$data[3,6] = 1.669069647789002
$data[3,7] = -0.8818392157554626
$data[4,0] = 400
$data[4,1] = "walkingToRunning"
$data[4,2] = 11.59089219358541
$data[4,3] = -11.35514022330562
$data[4,4] = -0.3842343388224623
$data[4,5] = 1.658749938011169
$data[4,6] = -1.585313200950623
$data[4,7] = -3.374025344848633
$data[5,0] = 500
$data[5,1] = "walkingToRunning"
$data[5,2] = 22.46794862013592
$data[5,3] = -12.58473284569012
$data[5,4] = 5.044853869274501
$data[5,5] = 2.416685581207275
$data[5,6] = -2.322476148605347
$data[5,7] = -5.666474342346191
$data[6,0] = 600
$data[6,1] = "walkingToRunning"
$data[6,2] = 3.344646702151348
$data[6,3] = -15.18862899802845
$data[6,4] = -5.453918620679495
$data[6,5] = 1.654488801956177
$data[6,6] = -1.7259281873703
$data[6,7] = -6.852646350860596
$data[7,0] = 700
$data[7,1] = "walkingToRunning"
$data[7,2] = -2.965254491602879
$data[7,3] = -16.62569452319625
$data[7,4] = 4.355243773150008
$data[7,5] = -3.016163110733032
$data[7,6] = -3.968842506408691
$data[7,7] = 3.157110929489136
$data[8,0] = 800
$data[8,1] = "walkingToRunning"
$data[8,2] = 3.611347477816996
$data[8,3] = -15.73287558132374
$data[8,4] = -0.8439290142623932
$data[8,5] = -5.865745544433594
$data[8,6] = 2.865894317626953
$data[8,7] = 8.381167411804199
$data[9,0] = 900
$data[9,1] = "walkingToRunning"
$data[9,2] = -2.01500372914868
$data[9,3] = -25.76902237164192
$data[9,4] = -10.39586599902989
$data[9,5] = -0.1431449055671692
$data[9,6] = 12.14008617401123
$data[9,7] = 1.953362107276917
$data[10,0] = 1000
$data[10,1] = "walkingToRunning"
$data[10,2] = 11.75571148212126
$data[10,3] = -16.33359923729541
$data[10,4] = 0.1073646178611298
$data[10,5] = 4.818324089050293
$data[10,6] = -6.334262371063232
$data[10,7] = -3.150319814682007
$data[11,0] = 1100
$data[11,1] = "walkingToRunning"
$data[11,2] = 39.60853858812318
$data[11,3] = -45.04104654887719
$data[11,4] = 25.44283379887685
$data[11,5] = 9.644504547119141
$data[11,6] = 3.98388934135437
$data[11,7] = 0.539756178855896
$data[12,0] = 1200
$data[12,1] = "walkingToRunning"
$data[12,2] = -21.39290519578826
$data[12,3] = -25.75901896854879
$data[12,4] = 6.399868327485425
$data[12,5] = 3.38767409324646
$data[12,6] = 3.137536764144897
$data[12,7] = -4.873384952545166
$data[13,0] = 1300
$data[13,1] = "walkingToRunning"
$data[13,2] = -34.85155020945177
$data[13,3] = -10.35765630959052
$data[13,4] = -4.708384362903586
$data[13,5] = -5.203151226043701
$data[13,6] = 0.9888983368873596
$data[13,7] = 2.413556337356567
$data[14,0] = 1400
$data[14,1] = "walkingToRunning"
$data[14,2] = -1.052037713090453
$data[14,3] = -12.47812661616762
$data[14,4] = -0.7732271075954033
$data[14,5] = -13.87813186645508
$data[14,6] = 0.5345630049705505
$data[14,7] = 0.425772875547409
$data[15,0] = 1500
$data[15,1] = "walkingToRunning"
$data[15,2] = 22.74343358553365
$data[15,3] = -12.66562087154951
$data[15,4] = 10.2766472348094
$data[15,5] = 5.411143779754639
$data[15,6] = -7.761183738708496
$data[15,7] = -2.545782327651977
$data[16,0] = 1600
$data[16,1] = "walkingToRunning"
$data[16,2] = 21.01473594981556
$data[16,3] = -12.33190507719506
$data[16,4] = 14.17957940750583
$data[16,5] = 2.374074935913086
$data[16,6] = -1.723265051841736
$data[16,7] = 1.63484799861908
$data[17,0] = 1700
$data[17,1] = "walkingToRunning"
$data[17,2] = 0.4384320727466777
$data[17,3] = -30.00794573225216
$data[17,4] = 5.117792562620172
$data[17,5] = 10.5457181930542
$data[17,6] = 8.709402084350586
$data[17,7] = -1.038965702056885
$data[18,0] = 1800
$data[18,1] = "walkingToRunning"
$data[18,2] = 1.731652502477468
$data[18,3] = -23.00902230076561
$data[18,4] = -15.75403659865687
$data[18,5] = -0.5905559659004211
$data[18,6] = -1.791441917419434
$data[18,7] = 0.2574611008167267
$data[19,0] = 1900
$data[19,1] = "walkingToRunning"
$data[19,2] = -8.11227344761256
$data[19,3] = 5.674290752974979
$data[19,4] = -28.99636612841348
$data[19,5] = -3.922703266143799
$data[19,6] = 3.348459005355835
$data[19,7] = 1.719536542892456
$data[20,0] = 2000
$data[20,1] = "walkingToRunning"
$data[20,2] = -38.96337933512132
$data[20,3] = -34.5937611280813
$data[20,4] = 3.551583487606536
$data[20,5] = -12.29035568237305
$data[20,6] = -13.32346248626709
$data[20,7] = -5.019326210021973
$data[21,0] = 2100
$data[21,1] = "walkingToRunning"
$data[21,2] = -22.12481265378448
$data[21,3] = -36.7449983704023
$data[21,4] = 14.93422636336839
$data[21,5] = 4.287290096282959
$data[21,6] = -8.83603572845459
$data[21,7] = -1.076782584190369
$data[22,0] = 2200
$data[22,1] = "walkingToRunning"
$data[22,2] = 1.603906876942109
$data[22,3] = -6.532943116137262
$data[22,4] = 5.020883938264563
$data[22,5] = 1.13570511341095
$data[22,6] = -4.037551879882812
$data[22,7] = 1.869738817214966
$data[23,0] = 2300
$data[23,1] = "walkingToRunning"
$data[23,2] = -0.676244735717717
$data[23,3] = -14.42456348125731
$data[23,4] = 5.205562151395338
$data[23,5] = 6.204765796661377
$data[23,6] = 5.688312530517578
$data[23,7] = -1.246159672737122
$data[24,0] = 2400
$data[24,1] = "walkingToRunning"
$data[24,2] = -8.344637724069424
$data[24,3] = -29.37389576646717
$data[24,4] = 4.180637551482669
$data[24,5] = -2.400972843170166
$data[24,6] = 2.221409320831299
$data[24,7] = 0.4785034656524658
$data[25,0] = 2500
$data[25,1] = "walkingToRunning"
$data[25,2] = -18.19668616345639
$data[25,3] = -0.7442300390208842
$data[25,4] = -9.97852061345025
$data[25,5] = -4.018577098846436
$data[25,6] = 3.744737386703491
$data[25,7] = 4.383230209350586
$data[26,0] = 2600
$data[26,1] = "walkingToRunning"
$data[26,2] = -38.98753061238127
$data[26,3] = -23.6091353738096
$data[26,4] = -9.86452339668957
$data[26,5] = -5.831124305725098
$data[26,6] = 10.37174797058106
$data[26,7] = -0.6757105588912964
$data[27,0] = 2700
$data[27,1] = "walkingToRunning"
$data[27,2] = -29.63654207760074
$data[27,3] = -38.7297830186655
$data[27,4] = 21.44472347208752
$data[27,5] = -0.996954381465912
$data[27,6] = 15.63681697845459
$data[27,7] = -8.430303573608398
$data[28,0] = 2800
$data[28,1] = "walkingToRunning"
$data[28,2] = -12.67720028352449
$data[28,3] = -6.331371307373074
$data[28,4] = 2.991470156336736
$data[28,5] = 5.841510772705078
$data[28,6] = -9.827264785766602
$data[28,7] = -3.761781692504883
$data[29,0] = 2900
$data[29,1] = "walkingToRunning"
$data[29,2] = 15.31067461092843
$data[29,3] = -20.30018748898473
$data[29,4] = 31.22906570886048
$data[29,5] = 1.826529026031494
$data[29,6] = 2.439255952835083
$data[29,7] = 2.031126499176025

$ws.Range("A2:H31").Value = $data
